# Data update for 5/10
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TN_AgeDaily")

# New date being appended: serial 43961 == 2020-05-10
$newDate = 43961

# Data rows: AgeBucket label, TOT_CASE_COUNT (C), DEATHS_TOT (G)
$data = @(
    @{ Age = "0-10";    C = 348;  G = 1  },
    @{ Age = "11-20";   C = 969;  G = 1  },
    @{ Age = "21-30";   C = 2939; G = 1  },
    @{ Age = "31-40";   C = 2886; G = 3  },
    @{ Age = "41-50";   C = 2632; G = 9  },
    @{ Age = "51-60";   C = 2402; G = 18 },
    @{ Age = "61-70";   C = 1462; G = 53 },
    @{ Age = "71-80";   C = 706;  G = 72 },
    @{ Age = "81+";     C = 445;  G = 85 },
    @{ Age = "Pending"; C = 196;  G = 0  }
)

$startRow = 532
$endRow = 541

# Touch the bottom-right new cell first so the sheet's used range already
# spans the new rows before any copy/paste happens (otherwise a PasteSpecial
# that extends the sheet leaves full-column refs like A:A blind to the new
# rows even though UsedRange reports the right extent).
$ws.Cells.Item($endRow, 9).Value = 0

# Copy formatting from the previous date's block (rows 522:531) down into the
# new block (rows 532:541) so number formats / fills match without minting
# new style entries.
$ws.Range("A522:I531").Copy()
$ws.Range("A532:I541").PasteSpecial(-4122)
$ws.Range("C522:C531").Copy()
$ws.Range("G532:G541").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$row = $startRow
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $newDate
    $ws.Cells.Item($row, 2).Value = $item.Age
    $ws.Cells.Item($row, 3).Value = $item.C
    $ws.Cells.Item($row, 7).Value = $item.G
    $row++
}

$row = $startRow
foreach ($item in $data) {
    $ws.Cells.Item($row, 4).Formula = "=C$row/SUMIF(A:A,A$row,C:C)"
    $ws.Cells.Item($row, 5).Formula = "=C$row-SUMIFS(C:C,A:A,A$row-1,B:B,B$row)"
    $ws.Cells.Item($row, 6).Formula = "=E$row/SUMIF(A:A,A$row,E:E)"
    $ws.Cells.Item($row, 8).Formula = "=G$row-SUMIFS(G:G,A:A,A$row-1,B:B,B$row)"
    $ws.Cells.Item($row, 9).Formula = "=G$row/SUMIF(A:A,A$row,G:G)"
    $row++
}

# Update sheet view to match new data extent (mirrors user scrolling to the new rows)
$ws.Application.ActiveWindow.ScrollRow = 528
$ws.Range("G532:G541").Select()
